$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the three new header cells (Wins, Losses, Ties) in AD1:AF1,
# copying the formatting (style) already used by the other header cells.
$ws.Range("AD1").Value2 = "Wins"
$ws.Range("AE1").Value2 = "Losses"
$ws.Range("AF1").Value2 = "Ties"

$ws.Range("AA1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill in the season record (Wins=71, Losses=91, Ties=0) for every data row.
$lastRow = 57
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value2 = 71
    $ws.Cells.Item($r, 31).Value2 = 91
    $ws.Cells.Item($r, 32).Value2 = 0
}
